$wb = $excel.ActiveWorkbook

# Delete columns C:F (vol min, vol max, Expected part, Left/right) on these 3 sheets
$sheetsToTrim = @("Clinical Structures", "opt structures", "couch_structures")
foreach ($name in $sheetsToTrim) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $range = $ws.Range("C1:F1048576")
    $range.Select()
    $range.Delete()
}

# Make couch_structures the active sheet/tab (was Doses before)
$ws4 = $wb.Worksheets.Item("couch_structures")
$ws4.Activate()
